$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values for rows 2-37, column G
$kValues = @(8,5,3,8,3,6,3,6,5,4,6,6,5,3,3,6,4,3,3,3,2,8,2,8,6,7,3,2,6,5,1,8,4,4,2,1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
